$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)
$ws.Activate()

# Row 9: replace "close browser" step with the new "Click on New Contact" step
$ws.Cells.Item(9, 1).Value = "Click on New Contact"
$ws.Cells.Item(9, 2).Value = "xpath"
$ws.Cells.Item(9, 3).Value = "//a[contains(text(),'New Contact')]"
$ws.Cells.Item(9, 4).Value = "click"
$ws.Cells.Item(9, 5).Value = "NA"

# Row 10: Enter Contact First Name
$ws.Cells.Item(10, 1).Value = "Enter Contact First Name"
$ws.Cells.Item(10, 2).Value = "xpath"
$ws.Cells.Item(10, 3).Value = "//input[@name='first_name' and @id='first_name']"
$ws.Cells.Item(10, 4).Value = "sendkeys"
$ws.Cells.Item(10, 5).Value = "Test"

# Row 11: Enter Contact Last Name
$ws.Cells.Item(11, 1).Value = "Enter Contact Last Name"
$ws.Cells.Item(11, 2).Value = "xpath"
$ws.Cells.Item(11, 3).Value = "//input[@name='surname' and @id='surname']"
$ws.Cells.Item(11, 4).Value = "sendkeys"
$ws.Cells.Item(11, 5).Value = "Test123"

# Row 12: Enter Contact Position
$ws.Cells.Item(12, 1).Value = "Enter Contact Position"
$ws.Cells.Item(12, 2).Value = "xpath"
$ws.Cells.Item(12, 3).Value = "//input[@name='company_position']"
$ws.Cells.Item(12, 4).Value = "sendkeys"
$ws.Cells.Item(12, 5).Value = "Manager"

# Row 13: Click on Save button
$ws.Cells.Item(13, 1).Value = "Click on Save button"
$ws.Cells.Item(13, 2).Value = "xpath"
$ws.Cells.Item(13, 3).Value = "//input[@type='submit' and @value='Save']"
$ws.Cells.Item(13, 4).Value = "click"
$ws.Cells.Item(13, 5).Value = "NA"

# Widen column C to fit the new, longer xpath values (closest achievable width to 43.58203125)
$ws.Columns.Item(3).ColumnWidth = 42.6

# Update the active selection to D8, matching the saved view state
$ws.Range("D8").Select() | Out-Null
